{"js": "// Item 2.1 (Declara\u00e7\u00e3o do Problema) of the Vision document: the text in the\n// \"O problema \u00e9\" table cell is corrected.\n//   old: \"an\u00fancios apresentados n\u00e3o s\u00e3o de interesse do usu\u00e1rio\"\n//   new: \"falta de informa\u00e7\u00f5es para proje\u00e7\u00e3o de anuncios\"\n\nconst oldText = \"an\u00fancios apresentados n\u00e3o s\u00e3o de interesse do usu\u00e1rio\";\nconst newText = \"falta de informa\u00e7\u00f5es para proje\u00e7\u00e3o de anuncios\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the matched range's text in place, preserving its formatting.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Item 2.1 (Declara\u00e7\u00e3o do Problema) of the Vision document: the text in the\n# \"O problema \u00e9\" table cell is corrected.\n#   old: \"an\u00fancios apresentados n\u00e3o s\u00e3o de interesse do usu\u00e1rio\"\n#   new: \"falta de informa\u00e7\u00f5es para proje\u00e7\u00e3o de anuncios\"\n\n$d = $word.ActiveDocument\n\n$oldText = \"an\u00fancios apresentados n\u00e3o s\u00e3o de interesse do usu\u00e1rio\"\n$newText = \"falta de informa\u00e7\u00f5es para proje\u00e7\u00e3o de anuncios\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
